# Applies the VIRGINIA_2017.xlsx cleanup edit:
#  1. Rename header columns to snake_case field names.
#  2. Title-case the Spanish connector words (de/del/la/las/el/los/y) that
#     appear inside state/municipality names in columns A and B.
#  3. Drop the trailing blank row + footnote/metadata rows (1151-1156),
#     shrinking the used range to A1:D1150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value2 = "mx_state"
$ws.Cells.Item(1, 2).Value2 = "mx_municipality"
$ws.Cells.Item(1, 3).Value2 = "n_matriculas"
$ws.Cells.Item(1, 4).Value2 = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B -------------------
function ProperConnectors($s) {
    $s = $s -replace '\bde\b', 'De'
    $s = $s -replace '\bdel\b', 'Del'
    $s = $s -replace '\bla\b', 'La'
    $s = $s -replace '\blas\b', 'Las'
    $s = $s -replace '\bel\b', 'El'
    $s = $s -replace '\blos\b', 'Los'
    $s = $s -replace '\by\b', 'Y'
    return $s
}

$lastDataRow = 1150
for ($r = 2; $r -le $lastDataRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value2 = ProperConnectors($a)
    }
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null) {
        $ws.Cells.Item($r, 2).Value2 = ProperConnectors($b)
    }
}

# --- 3. Tiny floating-point re-round of D434 (Zapopan) to match the
#        upstream recompute (58 / 6243, last-bit rounding) -------------
$ws.Cells.Item(434, 4).Value2 = 0.009290405253884353

# --- 4. Remove the trailing blank + footnote rows (1151-1156) -----------
$ws.Range("A1151:D1156").EntireRow.Delete()
